# Scene 24C edit ("write some new for stephen")
#
# The only authorable content change in this revision is a one-character
# trim: the final line of the scene,
#   "And being able to see that makes me pretty happy as well. "
# loses its trailing space, becoming
#   "And being able to see that makes me pretty happy as well."
#
# Locate it with Find/Replace over the whole document (there is exactly one
# occurrence) rather than hard-coding a paragraph index, so the script is
# resilient to the exact structure of the document.

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "And being able to see that makes me pretty happy as well. ",  # FindText
    $true,                                                          # MatchCase
    $false,                                                         # MatchWholeWord
    $false,                                                         # MatchWildcards
    $false,                                                         # MatchSoundsLike
    $false,                                                         # MatchAllWordForms
    $true,                                                          # Forward
    1,                                                              # Wrap: wdFindContinue
    $false,                                                         # Format
    "And being able to see that makes me pretty happy as well.",    # ReplaceWith
    2                                                                # Replace: wdReplaceAll
)
